$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 874.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 874.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 874.5
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").Value = -1100.5
$ws.Range("H33").Value = 342.36365
$ws.Range("I33").Value = 297.33334
$ws.Range("J33").Value = 545
$ws.Range("K33").Value = 297.33334
$ws.Range("L33").Value = 545
$ws.Range("M33").Value = -68.33334000000002
$ws.Range("N33").Value = -1003
$ws.Range("H40").Value = 6875
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 8750
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 8750
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = -9100
$ws.Range("H53").Value = 3139.2942
$ws.Range("I53").Value = 1444.3334
$ws.Range("J53").Value = 3502.5
$ws.Range("K53").Value = 1444.3334
$ws.Range("L53").Value = 3502.5
$ws.Range("M53").Value = -807.3334
$ws.Range("N53").Value = -4776.5
$ws.Range("H55").Value = 1934.625
$ws.Range("I55").Value = 487.58334
$ws.Range("J55").Value = 6275.75
$ws.Range("K55").Value = 487.58334
$ws.Range("L55").Value = 6275.75
$ws.Range("M55").Value = -273.58334
$ws.Range("N55").Value = -6703.75
$ws.Range("H80").Value = 2770.0938
$ws.Range("I80").Value = 1566.125
$ws.Range("J80").Value = 3974.0625
$ws.Range("K80").Value = 4698.375
$ws.Range("L80").Value = 11922.1875
$ws.Range("M80").Value = -3700.375
$ws.Range("N80").Value = -13918.1875
$ws.Range("H83").Value = 2770.0938
$ws.Range("I83").Value = 1566.125
$ws.Range("J83").Value = 3974.0625
$ws.Range("K83").Value = 14095.125
$ws.Range("L83").Value = 35766.5625
$ws.Range("M83").Value = -9103.125
$ws.Range("N83").Value = -45750.5625
$ws.Range("H86").Value = 4695.524
$ws.Range("I86").Value = 3180.7
$ws.Range("J86").Value = 6072.636
$ws.Range("K86").Value = 3180.7
$ws.Range("L86").Value = 6072.636
$ws.Range("M86").Value = -2057.7
$ws.Range("H89").Value = 4695.524
$ws.Range("I89").Value = 3180.7
$ws.Range("J89").Value = 6072.636
$ws.Range("K89").Value = 15903.5
$ws.Range("L89").Value = 30363.18
$ws.Range("M89").Value = -10287.5
$ws.Range("H101").Value = 1650.6666
$ws.Range("I101").Value = 1059.25
$ws.Range("J101").Value = 2833.5
$ws.Range("K101").Value = 3177.75
$ws.Range("L101").Value = 8500.5
$ws.Range("M101").Value = -1555.75
$ws.Range("N101").Value = -11744.5
$ws.Range("H103").Value = 1949.4
$ws.Range("I103").Value = 406.85715
$ws.Range("J103").Value = 3299.125
$ws.Range("K103").Value = 1220.57145
$ws.Range("L103").Value = 9897.375
$ws.Range("M103").Value = -634.5714499999999
$ws.Range("H106").Value = 2694.0557
$ws.Range("I106").Value = 1519.52
$ws.Range("J106").Value = 5363.4546
$ws.Range("K106").Value = 1519.52
$ws.Range("L106").Value = 5363.4546
$ws.Range("M106").Value = -888.52
$ws.Range("H107").Value = 259.85715
$ws.Range("I107").Value = 208.8
$ws.Range("J107").Value = 387.5
$ws.Range("K107").Value = 208.8
$ws.Range("L107").Value = 387.5
$ws.Range("M107").Value = 1711.2
$ws.Range("H112").Value = 1885.5238
$ws.Range("I112").Value = 1280
$ws.Range("J112").Value = 1915.8
$ws.Range("K112").Value = 3840
$ws.Range("L112").Value = 5747.4
$ws.Range("M112").Value = -2732
$ws.Range("N112").Value = -7963.4
$ws.Range("H137").Value = 11907266
$ws.Range("I137").Value = 31251432
$ws.Range("J137").Value = 3163.25
$ws.Range("K137").Value = 93754296
$ws.Range("L137").Value = 9489.75
$ws.Range("M137").Value = -93751746
$ws.Range("H141").Value = 2263.1428
$ws.Range("I141").Value = 2263.1428
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6789.428400000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1609.428400000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3487.8555
$ws.Range("I32").Value = 2646.0823
$ws.Range("J32").Value = 17798
$ws.Range("K32").Value = 2646.0823
$ws.Range("L32").Value = 17798
$ws.Range("M32").Value = -2359.0823
$ws.Range("N32").Value = -18372
$ws.Range("H51").Value = 40042
$ws.Range("I51").Value = 40042
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 40042
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -39286
$ws.Range("H61").Value = 5577
$ws.Range("I61").Value = 4364.25
$ws.Range("J61").Value = 7194
$ws.Range("K61").Value = 4364.25
$ws.Range("L61").Value = 7194
$ws.Range("M61").Value = -4152.25
$ws.Range("H74").Value = 13891744
$ws.Range("I74").Value = 23811902
$ws.Range("J74").Value = 3524.3
$ws.Range("K74").Value = 23811902
$ws.Range("L74").Value = 3524.3
$ws.Range("M74").Value = -23811028
$ws.Range("N74").Value = -5272.3
$ws.Range("H77").Value = 13891744
$ws.Range("I77").Value = 23811902
$ws.Range("J77").Value = 3524.3
$ws.Range("K77").Value = 119059510
$ws.Range("L77").Value = 17621.5
$ws.Range("M77").Value = -119055142
$ws.Range("N77").Value = -26357.5
$ws.Range("H110").Value = 7201.3
$ws.Range("I110").Value = 6074.2
$ws.Range("J110").Value = 10582.6
$ws.Range("K110").Value = 6074.2
$ws.Range("L110").Value = 10582.6
$ws.Range("M110").Value = -4029.2
$ws.Range("N110").Value = -14672.6
$ws.Range("H136").Value = 5577
$ws.Range("I136").Value = 4364.25
$ws.Range("J136").Value = 7194
$ws.Range("K136").Value = 13092.75
$ws.Range("L136").Value = 21582
$ws.Range("M136").Value = -10542.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5564.591
$ws.Range("I20").Value = 5010.353
$ws.Range("J20").Value = 7449
$ws.Range("K20").Value = 5010.353
$ws.Range("L20").Value = 7449
$ws.Range("M20").Value = -4763.353
$ws.Range("H86").Value = 3105.6333
$ws.Range("I86").Value = 3035.9375
$ws.Range("J86").Value = 3185.2856
$ws.Range("K86").Value = 3035.9375
$ws.Range("L86").Value = 3185.2856
$ws.Range("M86").Value = -1912.9375
$ws.Range("N86").Value = -5431.2856
$ws.Range("H89").Value = 3105.6333
$ws.Range("I89").Value = 3035.9375
$ws.Range("J89").Value = 3185.2856
$ws.Range("K89").Value = 15179.6875
$ws.Range("L89").Value = 15926.428
$ws.Range("M89").Value = -9563.6875
$ws.Range("N89").Value = -27158.428
$ws.Range("H132").Value = 61359.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 61359.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 61359.4
$ws.Range("N132").Value = -71479.4
$ws.Range("H134").Value = 1691.6786
$ws.Range("I134").Value = 1090.5
$ws.Range("J134").Value = 9507
$ws.Range("K134").Value = 3271.5
$ws.Range("L134").Value = 28521
$ws.Range("M134").Value = -736.5
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents() | Out-Null

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1290
$ws.Range("I16").Value = 685.1177
$ws.Range("J16").Value = 2759
$ws.Range("K16").Value = 685.1177
$ws.Range("L16").Value = 2759
$ws.Range("M16").Value = -398.1177
$ws.Range("H31").Value = 23723.084
$ws.Range("I31").Value = 1778.129
$ws.Range("J31").Value = 63740.35
$ws.Range("K31").Value = 1778.129
$ws.Range("L31").Value = 63740.35
$ws.Range("M31").Value = -1483.129
$ws.Range("N31").Value = -64330.35
$ws.Range("H34").Value = 23723.084
$ws.Range("I34").Value = 1778.129
$ws.Range("J34").Value = 63740.35
$ws.Range("K34").Value = 1778.129
$ws.Range("L34").Value = 63740.35
$ws.Range("M34").Value = -1576.129
$ws.Range("N34").Value = -64144.35
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 50000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 50000
$ws.Range("N48").Value = -50952
$ws.Range("H103").Value = 9886.333
$ws.Range("I103").Value = 9886.333
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 9886.333
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -8714.333
$ws.Range("H107").Value = 1430.2858
$ws.Range("I107").Value = 1195.875
$ws.Range("J107").Value = 1742.8334
$ws.Range("K107").Value = 1195.875
$ws.Range("L107").Value = 1742.8334
$ws.Range("M107").Value = 724.125
$ws.Range("N107").Value = -5582.8334
$ws.Range("H113").Value = 1290
$ws.Range("I113").Value = 685.1177
$ws.Range("J113").Value = 2759
$ws.Range("K113").Value = 685.1177
$ws.Range("L113").Value = 2759
$ws.Range("M113").Value = 1484.8823
$ws.Range("H132").Value = 4083.318
$ws.Range("I132").Value = 3332.5
$ws.Range("J132").Value = 6085.5
$ws.Range("K132").Value = 9997.5
$ws.Range("L132").Value = 18256.5
$ws.Range("M132").Value = -7467.5
$ws.Range("H134").Value = 3020.1155
$ws.Range("I134").Value = 2345.5789
$ws.Range("J134").Value = 4851
$ws.Range("K134").Value = 7036.736699999999
$ws.Range("L134").Value = 14553
$ws.Range("M134").Value = -4501.736699999999
$ws.Range("N134").Value = -19623

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27.666666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 27.666666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 82.999998
$ws.Range("N12").Value = -428.999998
$ws.Range("H39").Value = 6845.25
$ws.Range("I39").Value = 1248
$ws.Range("J39").Value = 8711
$ws.Range("K39").Value = 3744
$ws.Range("L39").Value = 26133
$ws.Range("M39").Value = -3450
$ws.Range("N39").Value = -26721
$ws.Range("H42").Value = 19002.5
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 19002.5
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 57007.5
$ws.Range("N42").Value = -58075.5
$ws.Range("H43").Value = 5490
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5490
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16470
$ws.Range("N43").Value = -16698
$ws.Range("H51").Value = 399.5
$ws.Range("I51").Value = 399.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1198.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -738.5
$ws.Range("N51").ClearContents() | Out-Null
$ws.Range("H81").Value = 2827.6365
$ws.Range("I81").Value = 1398.7778
$ws.Range("J81").Value = 9257.5
$ws.Range("K81").Value = 4196.3334
$ws.Range("L81").Value = 27772.5
$ws.Range("M81").Value = -3073.3334
$ws.Range("N81").Value = -30018.5
$ws.Range("H84").Value = 2827.6365
$ws.Range("I84").Value = 1398.7778
$ws.Range("J84").Value = 9257.5
$ws.Range("K84").Value = 12589.0002
$ws.Range("L84").Value = 83317.5
$ws.Range("M84").Value = -6973.0002
$ws.Range("N84").Value = -94549.5
$ws.Range("H129").Value = 6416259
$ws.Range("I129").Value = 2750
$ws.Range("J129").Value = 9266708
$ws.Range("K129").Value = 8250
$ws.Range("L129").Value = 27800124
$ws.Range("M129").Value = -3250
$ws.Range("N129").Value = -27810124
$ws.Range("H131").Value = 6521002
$ws.Range("I131").Value = 17858206
$ws.Range("J131").Value = 4631468
$ws.Range("K131").Value = 53574618
$ws.Range("L131").Value = 13894404
$ws.Range("M131").Value = -53569578
$ws.Range("N131").Value = -13904484

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 758.6667
$ws.Range("I2").Value = 173.25
$ws.Range("J2").Value = 1427.7142
$ws.Range("K2").Value = 173.25
$ws.Range("L2").Value = 1427.7142
$ws.Range("M2").Value = -60.25
$ws.Range("N2").Value = -1653.7142
$ws.Range("H33").Value = 30250
$ws.Range("I33").Value = 18000
$ws.Range("J33").Value = 34333.332
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 34333.332
$ws.Range("M33").Value = -17748
$ws.Range("N33").Value = -34837.332
$ws.Range("H36").Value = 14000
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 18500
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 18500
$ws.Range("M36").Value = -4515
$ws.Range("N36").Value = -19470
$ws.Range("H52").Value = 20000
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -19741
$ws.Range("H55").Value = 322598.8
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 322598.8
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 322598.8
$ws.Range("M55").ClearContents() | Out-Null
$ws.Range("N55").Value = -323252.8
$ws.Range("H102").Value = 2217
$ws.Range("I102").Value = 965.2692
$ws.Range("J102").Value = 5833.1113
$ws.Range("K102").Value = 965.2692
$ws.Range("L102").Value = 5833.1113
$ws.Range("M102").Value = 656.7308
$ws.Range("N102").Value = -9077.1113
$ws.Range("H122").Value = 5410.846
$ws.Range("I122").Value = 3464.8928
$ws.Range("J122").Value = 7681.125
$ws.Range("K122").Value = 10394.6784
$ws.Range("L122").Value = 23043.375
$ws.Range("M122").Value = -7944.678400000001
$ws.Range("N122").Value = -27943.375
$ws.Range("H132").Value = 2969.9678
$ws.Range("I132").Value = 2367.182
$ws.Range("J132").Value = 4443.4443
$ws.Range("K132").Value = 7101.545999999999
$ws.Range("L132").Value = 13330.3329
$ws.Range("M132").Value = -4571.545999999999
$ws.Range("N132").Value = -18390.3329

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3425.25
$ws.Range("I16").Value = 852
$ws.Range("J16").Value = 7027.8
$ws.Range("K16").Value = 852
$ws.Range("L16").Value = 7027.8
$ws.Range("M16").Value = -682
$ws.Range("N16").Value = -7367.8
$ws.Range("H40").Value = 6741.206
$ws.Range("I40").Value = 5545.7827
$ws.Range("J40").Value = 9240.728
$ws.Range("K40").Value = 5545.7827
$ws.Range("L40").Value = 9240.728
$ws.Range("M40").Value = -5409.7827
$ws.Range("H46").Value = 10545.546
$ws.Range("I46").Value = 6749.5
$ws.Range("J46").Value = 11389.111
$ws.Range("K46").Value = 6749.5
$ws.Range("L46").Value = 11389.111
$ws.Range("M46").Value = -6561.5
$ws.Range("N46").Value = -11765.111
$ws.Range("H50").Value = 35000
$ws.Range("I50").Value = 35000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 35000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -34363
$ws.Range("H53").Value = 15023
$ws.Range("I53").Value = 15023
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 15023
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -14505
$ws.Range("N53").ClearContents() | Out-Null
$ws.Range("H61").Value = 7635.4375
$ws.Range("I61").Value = 6665.5
$ws.Range("J61").Value = 9252
$ws.Range("K61").Value = 6665.5
$ws.Range("L61").Value = 9252
$ws.Range("M61").Value = -6463.5
$ws.Range("H82").Value = 7799.25
$ws.Range("I82").Value = 1564.6666
$ws.Range("J82").Value = 26503
$ws.Range("K82").Value = 1564.6666
$ws.Range("L82").Value = 26503
$ws.Range("M82").Value = -1203.6666
$ws.Range("N82").Value = -27225
$ws.Range("H85").Value = 7799.25
$ws.Range("I85").Value = 1564.6666
$ws.Range("J85").Value = 26503
$ws.Range("K85").Value = 1564.6666
$ws.Range("L85").Value = 26503
$ws.Range("M85").Value = -316.6666
$ws.Range("N85").Value = -28999
$ws.Range("H99").Value = 52000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 52000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 52000
$ws.Range("N99").Value = -57990
$ws.Range("H100").Value = 12225.104
$ws.Range("I100").Value = 9973.526
$ws.Range("J100").Value = 16503.1
$ws.Range("K100").Value = 9973.526
$ws.Range("L100").Value = 16503.1
$ws.Range("M100").Value = -9432.526
$ws.Range("N100").Value = -17585.1
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents() | Out-Null
$ws.Range("H113").Value = 7635.4375
$ws.Range("I113").Value = 6665.5
$ws.Range("J113").Value = 9252
$ws.Range("K113").Value = 6665.5
$ws.Range("L113").Value = 9252
$ws.Range("M113").Value = -4495.5
$ws.Range("H136").Value = 2618.746
$ws.Range("I136").Value = 1679.102
$ws.Range("J136").Value = 5907.5
$ws.Range("K136").Value = 5037.306
$ws.Range("L136").Value = 17722.5
$ws.Range("M136").Value = -2487.306
$ws.Range("N136").Value = -22822.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1187.7
$ws.Range("I107").Value = 1187.7
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3563.1
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1643.1
$ws.Range("H113").Value = 595.0323
$ws.Range("I113").Value = 268.875
$ws.Range("J113").Value = 1713.2858
$ws.Range("K113").Value = 806.625
$ws.Range("L113").Value = 5139.857400000001
$ws.Range("M113").Value = 1363.375
$ws.Range("H119").Value = 78947
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 78947
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 78947
$ws.Range("N119").Value = -88623
$ws.Range("H122").Value = 2961.9565
$ws.Range("I122").Value = 1880.375
$ws.Range("J122").Value = 5434.143
$ws.Range("K122").Value = 5641.125
$ws.Range("L122").Value = 16302.429
$ws.Range("M122").Value = -3191.125
$ws.Range("N122").Value = -21202.429
$ws.Range("H126").Value = 3668.7646
$ws.Range("I126").Value = 4818.8184
$ws.Range("J126").Value = 1560.3334
$ws.Range("K126").Value = 14456.4552
$ws.Range("L126").Value = 4681.0002
$ws.Range("M126").Value = -11986.4552
$ws.Range("N126").Value = -9621.0002
$ws.Range("H136").Value = 2904.2546
$ws.Range("I136").Value = 1908.75
$ws.Range("J136").Value = 6886.273
$ws.Range("K136").Value = 5726.25
$ws.Range("L136").Value = 20658.819
$ws.Range("M136").Value = -3176.25
$ws.Range("N136").Value = -25758.819
